$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Word wdHeaderFooterIndex constants
$wdHeaderFooterPrimary    = 1
$wdHeaderFooterFirstPage  = 2
$wdHeaderFooterEvenPages  = 3

# --- Enable / materialize the even-page header, first-page header and
#     first-page footer. Assigning content (even empty text + style) to
#     these mints the new header/footer parts and wires up the section's
#     headerReference/footerReference entries. The existing primary
#     (default) header keeps its original content and gets cloned to a
#     fresh part automatically, since it now needs to differ from the
#     freshly-created even-page header. ---

$hdrEven = $sec.Headers($wdHeaderFooterEvenPages)
$hdrEven.Range.Text = ""
$hdrEven.Range.Style = "Cabealho"

$hdrFirst = $sec.Headers($wdHeaderFooterFirstPage)
$hdrFirst.Range.Text = ""
$hdrFirst.Range.Style = "Cabealho"

$ftrFirst = $sec.Footers($wdHeaderFooterFirstPage)
$ftrFirst.Range.Text = ""
$ftrFirst.Range.Style = "Rodap"

# --- The primary/default header still reads "Casos de Testes do Projeto
#     Sistema de Rastreamento: Localizar Veículo" split across three
#     runs (the trailing ": " and "Localizar Veículo" are separate
#     runs). Re-running a Find/Replace over just the last run's text
#     (replacing it with itself) merges that trailing run into the
#     ": " run, leaving the heading's own run untouched. ---

$hdrPrimary = $sec.Headers($wdHeaderFooterPrimary)
$hdrPrimary.Range.Find.Execute("Localizar Veículo", $false, $false, $false, $false, $false, $true, 1, $false, "Localizar Veículo", 2)
